$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 7.878242263594639

# Row 3
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 7.700638116232206

# Row 4
$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 9.226618575922256
$ws.Range("D4").Value = 157.8057217802531
$ws.Range("E4").Value = 246.9852506941017
$ws.Range("G4").Value = 417.2004692788387
